# ============================================================
# LOB1285.docx edit: rotate body text among fixed paragraph
# slots (styles/positions unchanged; only text content moves).
#
# Implemented as a safe two-phase text rotation:
#   Phase 1: replace each distinct source text with a unique
#            placeholder token (this breaks the several
#            find/replace cycles created by the content move)
#   Phase 2: replace each placeholder token with its final
#            destination text.
#
# Find.Execute is used only to *locate* each exact range
# (Replace:=wdReplaceNone); the text is then assigned via
# Range.Text, which (unlike Find's own Replace) does not run
# AutoCorrect/AutoFormat smart-quote substitution, preserving
# the original characters exactly.
# ============================================================

$d = $word.ActiveDocument

function Set-ExactText($doc, [string]$oldText, [string]$newText) {
    $r = $doc.Content
    $found = $r.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $($oldText.Substring(0, [Math]::Min(60, $oldText.Length)))"
    }
    $r.Text = $newText
}

$txtA = @'
Orientar os estudantes no início de sua trajetória universitária no curso de graduação em Engenharia XX na EEL-USP de modo que o estudante seja capaz de a) identificar as oportunidades acadêmicas e as particularidades do seu curso; b) reconhecer, sob acompanhamento de um tutor, eventuais dificuldades ao longo do curso e compreender mecanismos para que estas sejam superadas, conduzindo o curso com o sucesso desejado; c) desenvolver habilidades técnicas e emocionais, ampliando as perspectivas de formação profissional por meio de atividades e encontros sistematizados.
'@
$txtB = @'
7455355 - Robson da Silva Rocha
'@
$txtDnorma = @'
Não se aplica
'@
$txtEcriterio = @'
Participação ativa nos encontros, apresentação de estudos/pesquisa e de trabalhos realizados durante a disciplina, colaboração e engajamento nas atividades da disciplina. A nota final é dada pela média ponderada das notas obtidas nas diversas atividades propostas.
'@
$txtFmetodo = @'
Atividades realizadas na forma de dinâmicas de grupos, utilização de vídeos, textos, roda de discussão e/ou elaboração de painéis. Participação em encontros de orientação promovidos pelo Programa de Tutoria Acadêmica e a realização de atividades propostas pelo tutor/monitor/mentor, incluindo trabalhos em equipe e estudos dirigidos.
'@
$txtG = @'
Apresentação dos programas e serviços oferecidos pela USP voltados aos estudantes e das oportunidades de realizar trabalhos extracurriculares. A dinâmica das aulas, ferramentas de interação. Desenvolvimento de atividades de grupo, com objetivo de desenvolver habilidades sócio-comportamentais através de colaboração em temas do curso relacionados à profissão escolhida. Áreas de atuação do curso de engenharia, competências e habilidades a serem desenvolvidas. Interdisciplinaridade e a relação entre as disciplinas e o conhecimento a ser aplicado. Planejamento de estudos. Formas de estudar e aprender.
'@
$txtH = @'
Os cursos de engenharia, respectivos projetos pedagógicos e seus componentes curriculares, incluindo TCC, estágio obrigatório, Projetos de Extensão Curricularizados, Atividades Acadêmicas Complementares e Atividades extracurriculares. Identificação e aderência do estudante com o curso e com a profissão escolhida. O curso superior, a transição adolescente/jovem adulto e os desafios nos projetos de vida do estudante no início da graduação. Relação entre as disciplinas e o conhecimento a ser aplicado. Competências e habilidades desenvolvidas no seu curso de engenharia. Dimensões acadêmicas, socioculturais e científicas. Diversidade e inclusão. Organização dos estudos.
'@
$txtI = @'
Guide students at the beginning of their university career in the XX Engineering undergraduate course at EEL-USP so that the student is able to a) identify the academic opportunities and particularities of their course; b) recognize, under the supervision of a tutor, any difficulties throughout the course and understand mechanisms for overcoming them, leading the course with the desired success; c) develop technical and emotional skills, expanding the perspectives of professional training through systematized activities and meetings.
'@
$txtJ = @'
Engineering courses, respective pedagogical projects and their curricular components, including TCC, mandatory internship, Curricular Extension Projects, Complementary Academic Activities and extracurricular activities. Identification and adherence of the student to the course and chosen profession. The higher education course, the adolescent/young adult transition and the challenges in the student's life projects at the beginning of graduation. Relationship between disciplines and the knowledge to be applied. Skills and abilities developed in your engineering course. Academic, sociocultural and scientific dimensions. Diversity and inclusion. Organization of studies.
'@
$bibSeg0 = @'
Bibliografia:
'@
$bibSeg1 = @'
A bibliografia será recomendada pelos docentes responsáveis e obtida na busca realizada pelos próprios alunos no início dos projetos. Seguem referências no tópico de mentoria: 
'@
$bibSeg2 = @'
[1] Peddy, S. The art of mentoring – Lead, follow and get out of the way. Houston: Bullion Books, 2001. 
'@
$bibSeg3 = @'
[2] Zachary, L. J. The Mentor’s Guide. San Francisco: Jossey-Bass Publishers, 2000. Pereira, A. Modelos de desenvolvimento do jovem adulto e promoção do bem-estar em estudantes do ensino superior. In: Programa de Monitorização e Tutorado: oito anos a promover a integração e o sucesso académico no IST. Lisboa: IST Press, 2011. p. 19-27. 
'@
$bibSeg4 = @'
[3] Mueller, S. Electronic mentoring as an example for the use of information and communications technology in engineering education. European Journal of Engineering Education, 2004. 
'@
$bibSeg5 = @'
[4] Kaul, S. Triangulated Mentorship of Engineering Students - Leveraging Peer Mentoring and Vertical Integration, Global Journal of Engineering Education, v. 21, p. 14-23,2019. 
'@
$bibSeg6 = @'
[5] Diretrizes Curriculares Nacionais para os cursos de graduação em Engenharia. Ministério da Educação. CNE/CES, 2019.
'@
$txtCbib = $bibSeg0 + "`v" + $bibSeg1 + "`v" + $bibSeg2 + "`v" + $bibSeg3 + "`v" + $bibSeg4 + "`v" + $bibSeg5 + "`v" + $bibSeg6

# ---------------- Phase 1: source text -> unique placeholder ----------------
Set-ExactText $d $txtA "PLACEHOLDERONE"
Set-ExactText $d $txtI "PLACEHOLDERTWO"
Set-ExactText $d $txtB "PLACEHOLDERTHREE"
Set-ExactText $d $txtH "PLACEHOLDERFOUR"
Set-ExactText $d $txtJ "PLACEHOLDERFIVE"
Set-ExactText $d $txtG "PLACEHOLDERSIX"
Set-ExactText $d $txtFmetodo "PLACEHOLDERSEVEN"
Set-ExactText $d $txtEcriterio "PLACEHOLDEREIGHT"
Set-ExactText $d $txtDnorma "PLACEHOLDERNINE"
Set-ExactText $d $txtCbib "PLACEHOLDERTEN"

# ---------------- Phase 2: placeholder -> final text ----------------
Set-ExactText $d "PLACEHOLDERONE" $txtH
Set-ExactText $d "PLACEHOLDERTWO" $txtJ
Set-ExactText $d "PLACEHOLDERTHREE" $txtA
Set-ExactText $d "PLACEHOLDERFOUR" $txtG
Set-ExactText $d "PLACEHOLDERFIVE" $txtI
Set-ExactText $d "PLACEHOLDERSIX" $txtFmetodo
Set-ExactText $d "PLACEHOLDERSEVEN" $txtEcriterio
Set-ExactText $d "PLACEHOLDEREIGHT" $txtDnorma
Set-ExactText $d "PLACEHOLDERNINE" $txtCbib
Set-ExactText $d "PLACEHOLDERTEN" $txtB

Write-Host "Done."
